# Applies the "Statis metricsdb added locally" edit to Sheet1:
#  - expands the comparison table from 2 tickers x 4 metrics to
#    6 tickers x 11 metrics, moving the ticker symbol out of column A
#    (which becomes a numeric 0-filled index column) into column B.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replicate the existing bold/centered/bordered "index" style (currently
#     on A2:A3) down onto the newly added rows A4:A7 before touching values ---
$ws.Range("A3").Copy()
$ws.Range("A4:A7").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Likewise replicate the navy/bold/bordered header style (currently on
#     E1) onto the new trailing header cells F1:L1 ---
$ws.Range("E1").Copy()
$ws.Range("F1:L1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Header row (row 1): renamed + new trailing metric columns ---
$ws.Cells.Item(1,2).Value = "ticker"
$ws.Cells.Item(1,3).Value = "Debt growth"
$ws.Cells.Item(1,4).Value = "Net income available to common shareholders"
$ws.Cells.Item(1,5).Value = "Other working capital"
$ws.Cells.Item(1,6).Value = "Net income per EBT"
$ws.Cells.Item(1,7).Value = "Dividend yield"
$ws.Cells.Item(1,8).Value = "Operating cash flow sales ratio"
$ws.Cells.Item(1,9).Value = "Current ratio"
$ws.Cells.Item(1,10).Value = "Total liabilities"
$ws.Cells.Item(1,11).Value = "Debt repayment"
$ws.Cells.Item(1,12).Value = "Net cash used provided by (used for) financing activities"

# --- Column widths for the newly added columns E..K (match existing A..D) ---
$ws.Columns.Item(5).ColumnWidth = 12.140625
$ws.Columns.Item(6).ColumnWidth = 12.140625
$ws.Columns.Item(7).ColumnWidth = 12.140625
$ws.Columns.Item(8).ColumnWidth = 12.140625
$ws.Columns.Item(9).ColumnWidth = 12.140625
$ws.Columns.Item(10).ColumnWidth = 12.140625
$ws.Columns.Item(11).ColumnWidth = 12.140625

# --- Data rows: col A = numeric placeholder index (0), col B = ticker, ---
# --- cols C..L = the metric values ---
# row 2: AAPL
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "AAPL"
$ws.Cells.Item(2,3).Value = -0.05621795375732642
$ws.Cells.Item(2,4).Value = 94760
$ws.Cells.Item(2,5).Value = 100987
$ws.Cells.Item(2,6).Value = 0.8528082577196314
$ws.Cells.Item(2,7).Value = 0.005573960673721321
$ws.Cells.Item(2,8).Value = 0.2884094081427658
$ws.Cells.Item(2,9).Value = 0.9880116717592975
$ws.Cells.Item(2,10).Value = 290437
$ws.Cells.Item(2,11).Value = -13944000000
$ws.Cells.Item(2,12).Value = -112129000000

# row 3: NVDA
$ws.Cells.Item(3,1).Value = 0
$ws.Cells.Item(3,2).Value = "NVDA"
$ws.Cells.Item(3,3).Value = -0.006
$ws.Cells.Item(3,4).Value = 10326
$ws.Cells.Item(3,5).Value = 10323
$ws.Cells.Item(3,6).Value = 1.044726142071275
$ws.Cells.Item(3,7).Value = 0.0008351537797192516
$ws.Cells.Item(3,8).Value = 0.2091273077778601
$ws.Cells.Item(3,9).Value = 3.515617857687033
$ws.Cells.Item(3,10).Value = 19081
$ws.Cells.Item(3,11).Value = -1250000000
$ws.Cells.Item(3,12).Value = -10888000000

# row 4: AMZN
$ws.Cells.Item(4,1).Value = 0
$ws.Cells.Item(4,2).Value = "AMZN"
$ws.Cells.Item(4,3).Value = 0.04906924009914043
$ws.Cells.Item(4,4).Value = 20079
$ws.Cells.Item(4,5).Value = 16921
$ws.Cells.Item(4,6).Value = 0.4583263175618791
$ws.Cells.Item(4,7).Value = 0
$ws.Cells.Item(4,8).Value = 0.09096020685509054
$ws.Cells.Item(4,9).Value = 0.9446435811136924
$ws.Cells.Item(4,10).Value = 316632
$ws.Cells.Item(4,11).Value = -45272000000
$ws.Cells.Item(4,12).Value = -9047000000

# row 5: GOOGL
$ws.Cells.Item(5,1).Value = 0
$ws.Cells.Item(5,2).Value = "GOOGL"
$ws.Cells.Item(5,3).Value = 0.0108339632149156
$ws.Cells.Item(5,4).Value = 66732
$ws.Cells.Item(5,5).Value = 77618
$ws.Cells.Item(5,6).Value = 0.8407918349035441
$ws.Cells.Item(5,7).Value = 0
$ws.Cells.Item(5,8).Value = 0.3234913518788273
$ws.Cells.Item(5,9).Value = 2.377994227994228
$ws.Cells.Item(5,10).Value = 109120
$ws.Cells.Item(5,11).Value = -18339000000
$ws.Cells.Item(5,12).Value = -70414000000

# row 6: MSFT
$ws.Cells.Item(6,1).Value = 0
$ws.Cells.Item(6,2).Value = "MSFT"
$ws.Cells.Item(6,3).Value = -0.04207412477997262
$ws.Cells.Item(6,4).Value = 77096
$ws.Cells.Item(6,5).Value = 63226
$ws.Cells.Item(6,6).Value = 0.8102137474667174
$ws.Cells.Item(6,7).Value = 0.007777998514954823
$ws.Cells.Item(6,8).Value = 0.4132883467427978
$ws.Cells.Item(6,9).Value = 1.76916725076573
$ws.Cells.Item(6,10).Value = 205753
$ws.Cells.Item(6,11).Value = -3250000000
$ws.Cells.Item(6,12).Value = -18291000000

# row 7: GOOG
$ws.Cells.Item(7,1).Value = 0
$ws.Cells.Item(7,2).Value = "GOOG"
$ws.Cells.Item(7,3).Value = 0.0108339632149156
$ws.Cells.Item(7,4).Value = 66732
$ws.Cells.Item(7,5).Value = 77618
$ws.Cells.Item(7,6).Value = 0.8407918349035441
$ws.Cells.Item(7,7).Value = 0
$ws.Cells.Item(7,8).Value = 0.3234913518788273
$ws.Cells.Item(7,9).Value = 2.377994227994228
$ws.Cells.Item(7,10).Value = 109120
$ws.Cells.Item(7,11).Value = -18339000000
$ws.Cells.Item(7,12).Value = -70414000000

# --- C2/C3 previously carried a 0.00% number format; the refreshed dump
#     uses the default General format there, so reset those two cells ---
$ws.Cells.Item(2,3).Style = "Normal"
$ws.Cells.Item(3,3).Style = "Normal"
